# Updated cryptos list on Mon May 13 13:48:48 UTC 2024 with GitHub Actions
#
# Refreshes the coinranking price/volume snapshot on Sheet1 (columns
# B:E, rows 2-51): latest Price (D) and Volume(1h) (E) for every coin,
# plus two rows (24/25) that swapped rank position (Litecoin <-> RenderToken)
# so their Coin/Link/Price/Volume all move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new display text, taken straight from the refreshed feed.
$updates = @(
    @{ Cell = 'D2'; Value = '62.696.49' }
    @{ Cell = 'E2'; Value = '  +2.51%  ' }
    @{ Cell = 'D3'; Value = '2.967.15' }
    @{ Cell = 'E3'; Value = '  +1.15%  ' }
    @{ Cell = 'E4'; Value = '  -0.01%  ' }
    @{ Cell = 'D5'; Value = '594.62' }
    @{ Cell = 'E5'; Value = '  +0.20%  ' }
    @{ Cell = 'D6'; Value = '145.38' }
    @{ Cell = 'E6'; Value = '  +0.15%  ' }
    @{ Cell = 'E7'; Value = '  +0.07%  ' }
    @{ Cell = 'D8'; Value = '2.964.44' }
    @{ Cell = 'E8'; Value = '  +1.08%  ' }
    @{ Cell = 'E9'; Value = '  +0.09%  ' }
    @{ Cell = 'E10'; Value = '  +3.46%  ' }
    @{ Cell = 'D11'; Value = '0.147' }
    @{ Cell = 'E11'; Value = '  +2.70%  ' }
    @{ Cell = 'D12'; Value = '0.446' }
    @{ Cell = 'E12'; Value = '  +0.70%  ' }
    @{ Cell = 'E13'; Value = '  +5.29%  ' }
    @{ Cell = 'D14'; Value = '33.08' }
    @{ Cell = 'E14'; Value = '  -2.10%  ' }
    @{ Cell = 'E15'; Value = '  -0.49%  ' }
    @{ Cell = 'D16'; Value = '3.457.17' }
    @{ Cell = 'E16'; Value = '  +1.09%  ' }
    @{ Cell = 'D17'; Value = '62.585.08' }
    @{ Cell = 'E17'; Value = '  +2.45%  ' }
    @{ Cell = 'E18'; Value = '  -0.49%  ' }
    @{ Cell = 'D19'; Value = '2.947.02' }
    @{ Cell = 'E19'; Value = '  +0.50%  ' }
    @{ Cell = 'D20'; Value = '441.50' }
    @{ Cell = 'E20'; Value = '  +1.70%  ' }
    @{ Cell = 'D21'; Value = '13.45' }
    @{ Cell = 'E21'; Value = '  -0.46%  ' }
    @{ Cell = 'E22'; Value = '  -1.49%  ' }
    @{ Cell = 'D23'; Value = '7.06' }
    @{ Cell = 'E23'; Value = '  -0.89%  ' }
    @{ Cell = 'B24'; Value = 'RenderToken' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D24'; Value = '11.28' }
    @{ Cell = 'E24'; Value = '  +1.67%  ' }
    @{ Cell = 'B25'; Value = 'Litecoin' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D25'; Value = '81.69' }
    @{ Cell = 'E25'; Value = '  -0.06%  ' }
    @{ Cell = 'E26'; Value = '  +0.23%  ' }
    @{ Cell = 'D27'; Value = '2.12' }
    @{ Cell = 'E27'; Value = '  -3.92%  ' }
    @{ Cell = 'E28'; Value = '  +0.01%  ' }
    @{ Cell = 'D29'; Value = '7.25' }
    @{ Cell = 'E29'; Value = '  +3.60%  ' }
    @{ Cell = 'E30'; Value = '  -0.19%  ' }
    @{ Cell = 'E31'; Value = '  -4.70%  ' }
    @{ Cell = 'D32'; Value = '0.0₃0952' }
    @{ Cell = 'E32'; Value = '  +9.33%  ' }
    @{ Cell = 'E33'; Value = '  -0.15%  ' }
    @{ Cell = 'E34'; Value = '  -1.09%  ' }
    @{ Cell = 'E35'; Value = '  -0.05%  ' }
    @{ Cell = 'D36'; Value = '0.993' }
    @{ Cell = 'E36'; Value = '  -2.13%  ' }
    @{ Cell = 'E37'; Value = '  -0.45%  ' }
    @{ Cell = 'E38'; Value = '  +1.67%  ' }
    @{ Cell = 'D39'; Value = '2.06' }
    @{ Cell = 'E39'; Value = '  +3.12%  ' }
    @{ Cell = 'D40'; Value = '49.46' }
    @{ Cell = 'E40'; Value = '  -0.98%  ' }
    @{ Cell = 'E41'; Value = '  -0.90%  ' }
    @{ Cell = 'E42'; Value = '  -4.82%  ' }
    @{ Cell = 'E43'; Value = '  -0.73%  ' }
    @{ Cell = 'D44'; Value = '40.20' }
    @{ Cell = 'E44'; Value = '  -4.51%  ' }
    @{ Cell = 'D45'; Value = '2.743.07' }
    @{ Cell = 'E45'; Value = '  +1.26%  ' }
    @{ Cell = 'D46'; Value = '135.33' }
    @{ Cell = 'E46'; Value = '  +1.11%  ' }
    @{ Cell = 'E47'; Value = '  -1.95%  ' }
    @{ Cell = 'D48'; Value = '362.34' }
    @{ Cell = 'E48'; Value = '  -3.35%  ' }
    @{ Cell = 'E49'; Value = '  +0.03%  ' }
    @{ Cell = 'D50'; Value = '23.04' }
    @{ Cell = 'E50'; Value = '  -3.77%  ' }
    @{ Cell = 'E51'; Value = '  -0.66%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $value = $u.Value

    # Column D holds price strings such as "594.62" or "62.696.49" that
    # read as numeric (or as a date, for the two-dot "thousand.thousand"
    # style prices). Writing them straight through `.Value` lets Excel's
    # usual type-sniffing coerce them into a number/date and silently
    # normalise the text (e.g. "441.50" -> 441.5), which would not match
    # the original inline-string cell. A leading apostrophe forces plain
    # text and preserves the exact digits/trailing zeros.
    $looksNumeric = $value -match '^-?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $cell.Value = "'" + $value
        # The apostrophe-prefix entry tags the cell with Excel's
        # "number stored as text" quote-prefix style. Re-stamp the style
        # from an already-plain text cell (B2, the "Coin" header's first
        # data cell) so the cell keeps its original, unstyled look.
        $cell.Style = $ws.Range("B2").Style
    } else {
        $cell.Value = $value
    }
}
